$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Date column (F) for all data rows to the new timestamp
$newDate = "Feb 12, 2022 (04:14:47 EST)"
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 6).Value = $newDate
}

# Fix B7 from "BOMB" to "BOM"
$ws.Cells.Item(7, 2).Value = "BOM"
